# Add two new columns, "I0" (I) and "IF" (J), to the results table,
# mirroring the header formatting already used by the neighbouring "IP"
# column (bold header, border, centered) and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the existing header style (bold font + border + center alignment)
# from H1 onto the two new header cells before writing their text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I2:J12 values, row-aligned with the existing data (rows 2-12).
$values = @(
    @(8, 8),
    @(6, 7),
    @(2, 4),
    @(6, 8),
    @(6, 6),
    @(6, 9),
    @(2, 5),
    @(6, 6),
    @(4, 5),
    @(3, 4),
    @(1, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]   # column I
    $ws.Cells.Item($row, 10).Value = $values[$i][1]  # column J
}
